$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L3").NumberFormat = "General"
$ws.Range("L3").Value = "FALSE"
